# Apply the edits described by the diff:
# 1. Change D4 from 41 to 45.
# 2. Clear out the data in row 6 (values removed, but keep the per-column styles).
# 3. Extend the used range down to row 7 by touching A7:F7, then clear its contents
#    so the new row 7 ends up fully empty (same shape as the cleared row 6), which
#    also grows the sheet dimension from A1:F6 to A1:F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the counter value in D4.
$ws.Range("D4").Value = 45

# Wipe out the GGG / 835DF613 / 104|101 / TRUE row (row 6), leaving empty, styled cells.
$ws.Range("A6:F6").ClearContents()

# Materialize a new empty row 7 (same styling pattern as the other rows) and extend
# the sheet's used range/dimension to A1:F7.
$ws.Range("A7:F7").Value = 0
$ws.Range("A7:F7").ClearContents()
